$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing obstacle-loop-check table (rows 13-22, cols N:O) one
# column to the right (N->O, O->P) to make room for a new column that will
# hold the candidate obstacle coordinate, mirroring "copying the grid object
# to add a potential obstacle". Range.Insert() in this runtime shifts whole
# columns rather than just the selected rows, so move the block manually:
# copy the rightmost column first (O->P), then the left one (N->O), each as
# a separate values pass + formats pass so per-cell styles move correctly.
$ws.Range("O13:O22").Copy()
$ws.Range("P13:P22").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("O13:O22").Copy()
$ws.Range("P13:P22").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("N13:N22").Copy()
$ws.Range("O13:O22").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("N13:N22").Copy()
$ws.Range("O13:O22").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# N13:N16 have no new content - drop the leftover copies of the old column
$ws.Range("N13:N16").Clear()

# --- New header row 1: column index labels (0-9) above the grid, in N1:W1
for ($i = 0; $i -le 9; $i++) {
    $col = 14 + $i  # N=14 .. W=23
    $ws.Cells.Item(1, $col).Value = $i
}

# --- New column M: row index labels (0-9) beside the grid, in M2:M11
for ($i = 0; $i -le 9; $i++) {
    $row = 2 + $i  # rows 2..11
    $ws.Cells.Item($row, 13).Value = $i  # M=13
}

# --- New column N (rows 17-22): the candidate obstacle coordinate being
# tested for each loop check entry
$ws.Range("N17").Value = "3,6"
$ws.Range("N18").Value = "7,6"
$ws.Range("N19").Value = "3,8"
$ws.Range("N20").Value = "1,8"
$ws.Range("N21").Value = "7,7"
$ws.Range("N22").Value = "7,9"

# --- Update the active selection to reflect where editing left off
$ws.Range("N23").Select()

Write-Host "done"
